$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 30065.885
$ws.Range("I11").Value = 30065.885
$ws.Range("K11").Value = 30065.885
$ws.Range("M11").Value = -29925.885

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 145.85715
$ws.Range("I33").Value = 154.33333
$ws.Range("K33").Value = 154.33333
$ws.Range("M33").Value = 74.66667000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4764212
$ws.Range("I137").Value = 7409027.5
$ws.Range("J137").Value = 3543.8667
$ws.Range("K137").Value = 22227082.5
$ws.Range("L137").Value = 10631.6001
$ws.Range("M137").Value = -22224532.5
$ws.Range("N137").Value = -15731.6001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5187.717
$ws.Range("J138").Value = 5902.278
$ws.Range("L138").Value = 17706.834
$ws.Range("N138").Value = -27986.834

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 1551.2
$ws.Range("I3").Value = 2250
$ws.Range("K3").Value = 2250
$ws.Range("M3").Value = -2135

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 939943.8
$ws.Range("I32").Value = 1203558.5
$ws.Range("K32").Value = 1203558.5
$ws.Range("M32").Value = -1203271.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4903.3335
$ws.Range("I63").Value = 2938.75
$ws.Range("K63").Value = 2938.75
$ws.Range("M63").Value = -2252.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 4903.3335
$ws.Range("I66").Value = 2938.75
$ws.Range("K66").Value = 14693.75
$ws.Range("M66").Value = -11261.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2302.52
$ws.Range("I88").Value = 2429.1667
$ws.Range("J88").Value = 2185.6155
$ws.Range("K88").Value = 2429.1667
$ws.Range("L88").Value = 2185.6155
$ws.Range("M88").Value = -2023.1667
$ws.Range("N88").Value = -2997.6155

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2302.52
$ws.Range("I91").Value = 2429.1667
$ws.Range("J91").Value = 2185.6155
$ws.Range("K91").Value = 2429.1667
$ws.Range("L91").Value = 2185.6155
$ws.Range("M91").Value = -1025.1667
$ws.Range("N91").Value = -4993.6155

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 910.4706
$ws.Range("I97").Value = 910.4706
$ws.Range("K97").Value = 910.4706
$ws.Range("M97").Value = -414.4706

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1795.4615
$ws.Range("I102").Value = 1795.4615
$ws.Range("K102").Value = 1795.4615
$ws.Range("M102").Value = -173.4614999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 6500
$ws.Range("I25").Value = 6500
$ws.Range("K25").Value = 6500
$ws.Range("M25").Value = -6265

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4233.143
$ws.Range("I105").Value = 4194.8335
$ws.Range("J105").Value = 4261.875
$ws.Range("K105").Value = 4194.8335
$ws.Range("L105").Value = 4261.875
$ws.Range("M105").Value = -2447.8335
$ws.Range("N105").Value = -7755.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 1500
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 1500
$ws.Range("K14").Value = 0
$ws.Range("N14").Value = -1840
$ws.Range("M14").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 57647.11
$ws.Range("I16").Value = 2015
$ws.Range("J16").Value = 202290.6
$ws.Range("K16").Value = 2015
$ws.Range("L16").Value = 202290.6
$ws.Range("M16").Value = -1728
$ws.Range("N16").Value = -202864.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1252604.1
$ws.Range("I31").Value = 1966073.4
$ws.Range("J31").Value = 4032.8333
$ws.Range("K31").Value = 1966073.4
$ws.Range("L31").Value = 4032.8333
$ws.Range("M31").Value = -1965778.4
$ws.Range("N31").Value = -4622.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1252604.1
$ws.Range("I34").Value = 1966073.4
$ws.Range("J34").Value = 4032.8333
$ws.Range("K34").Value = 1966073.4
$ws.Range("L34").Value = 4032.8333
$ws.Range("M34").Value = -1965871.4
$ws.Range("N34").Value = -4436.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("N70").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 57647.11
$ws.Range("I113").Value = 2015
$ws.Range("J113").Value = 202290.6
$ws.Range("K113").Value = 2015
$ws.Range("L113").Value = 202290.6
$ws.Range("M113").Value = 155
$ws.Range("N113").Value = -206630.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H116").Value = 44998
$ws.Range("J116").Value = 44998
$ws.Range("L116").Value = 44998
$ws.Range("N116").Value = -54176

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 3354.182
$ws.Range("J107").Value = 4924.0713
$ws.Range("L107").Value = 14772.2139
$ws.Range("N107").Value = -18612.2139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1276.3572
$ws.Range("I114").Value = 173.33333
$ws.Range("J114").Value = 2103.625
$ws.Range("K114").Value = 519.99999
$ws.Range("L114").Value = 6310.875
$ws.Range("M114").Value = 2734.00001
$ws.Range("N114").Value = -12818.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 6274.0625
$ws.Range("J127").Value = 6274.0625
$ws.Range("L127").Value = 18822.1875
$ws.Range("N127").Value = -28742.1875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 456733.1
$ws.Range("I129").Value = 1251988.5
$ws.Range("J129").Value = 2301.4285
$ws.Range("K129").Value = 3755965.5
$ws.Range("L129").Value = 6904.2855
$ws.Range("M129").Value = -3750965.5
$ws.Range("N129").Value = -16904.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 6365.909
$ws.Range("I137").Value = 2407.7
$ws.Range("J137").Value = 9664.416999999999
$ws.Range("K137").Value = 7223.099999999999
$ws.Range("L137").Value = 28993.251
$ws.Range("M137").Value = -2123.099999999999
$ws.Range("N137").Value = -39193.251

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3034.2307
$ws.Range("I80").Value = 1978
$ws.Range("K80").Value = 1978
$ws.Range("M80").Value = -980

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3034.2307
$ws.Range("I83").Value = 1978
$ws.Range("K83").Value = 9890
$ws.Range("M83").Value = -4898

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2445.875
$ws.Range("I102").Value = 2295.2856
$ws.Range("J102").Value = 3500
$ws.Range("K102").Value = 2295.2856
$ws.Range("L102").Value = 3500
$ws.Range("M102").Value = -673.2856000000002
$ws.Range("N102").Value = -6744

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 59123.39
$ws.Range("I122").Value = 65751.31
$ws.Range("J122").Value = 6100
$ws.Range("K122").Value = 197253.93
$ws.Range("L122").Value = 18300
$ws.Range("M122").Value = -194803.93
$ws.Range("N122").Value = -23200

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2853.3635
$ws.Range("I7").Value = 2049
$ws.Range("K7").Value = 2049
$ws.Range("M7").Value = -1937

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3454.4736
$ws.Range("I22").Value = 2652.8333
$ws.Range("K22").Value = 2652.8333
$ws.Range("M22").Value = -2357.8333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 3454.4736
$ws.Range("I27").Value = 2652.8333
$ws.Range("K27").Value = 2652.8333
$ws.Range("M27").Value = -2545.8333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4042.3333
$ws.Range("I40").Value = 3510.4614
$ws.Range("K40").Value = 3510.4614
$ws.Range("M40").Value = -3374.4614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1471.2727
$ws.Range("I82").Value = 1471.2727
$ws.Range("K82").Value = 1471.2727
$ws.Range("M82").Value = -1110.2727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1471.2727
$ws.Range("I85").Value = 1471.2727
$ws.Range("K85").Value = 1471.2727
$ws.Range("M85").Value = -223.2727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5227.846
$ws.Range("I122").Value = 4496.4443
$ws.Range("K122").Value = 13489.3329
$ws.Range("M122").Value = -11039.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2853.3635
$ws.Range("I126").Value = 2049
$ws.Range("K126").Value = 6147
$ws.Range("M126").Value = -3677

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2087032.5
$ws.Range("I132").Value = 3707519
$ws.Range("J132").Value = 3550.1428
$ws.Range("K132").Value = 11122557
$ws.Range("L132").Value = 10650.4284
$ws.Range("M132").Value = -11120027
$ws.Range("N132").Value = -15710.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 500
$ws.Range("J6").Value = 500
$ws.Range("L6").Value = 500
$ws.Range("N6").Value = -730

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 50000
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 66116.81
$ws.Range("I81").Value = 3037.4614
$ws.Range("J81").Value = 339460.66
$ws.Range("K81").Value = 6074.9228
$ws.Range("L81").Value = 678921.3199999999
$ws.Range("M81").Value = -5013.9228
$ws.Range("N81").Value = -681043.3199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 66116.81
$ws.Range("I84").Value = 3037.4614
$ws.Range("J84").Value = 339460.66
$ws.Range("K84").Value = 30374.614
$ws.Range("L84").Value = 3394606.6
$ws.Range("M84").Value = -25070.614
$ws.Range("N84").Value = -3405214.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3331.1943
$ws.Range("I122").Value = 1970.4062
$ws.Range("K122").Value = 5911.2186
$ws.Range("M122").Value = -3461.2186

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4177.2
$ws.Range("I126").Value = 2575.7273
$ws.Range("K126").Value = 7727.1819
$ws.Range("M126").Value = -5257.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 20837796
$ws.Range("I132").Value = 41671236
$ws.Range("J132").Value = 4354.75
$ws.Range("K132").Value = 125013708
$ws.Range("L132").Value = 13064.25
$ws.Range("M132").Value = -125011178
$ws.Range("N132").Value = -18124.25
